# "feat: add 2022-Q3 data"
#
# The workbook tracks BAYN (Bayer AG) holdings per quarter. A new quarter
# "2022-Q3" is being added:
#   - a new per-quarter worksheet "2022-Q3" is inserted right before the
#     existing "2022-Q2" worksheet, holding the new quarter's fund data;
#   - the "总计" (totals) summary worksheet gets a new data row for
#     "2022-Q3" inserted at the top of its data (row 2), with every
#     following quarter's row pushed down by one.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: add "2022-Q3" as the newest
#    entry and push the older quarters down one row.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("D2").Value = 0.24

$totals.Range("B3").Value = "2022-Q2"
$totals.Range("D3").Value = 0.29

$totals.Range("B4").Value = "2022-Q1"
$totals.Range("D4").Value = 0.29

$totals.Range("B5").Value = "2021-Q2"
$totals.Range("D5").Value = 0.32

$totals.Range("B6").Value = "2021-Q1"
$totals.Range("D6").Value = 0.38

# New row 7 ("2020-Q4") - copy A6's number format down to A7 first so the
# new row matches the look of the existing index column, then fill values.
$totals.Range("A6").Copy()
$totals.Range("A7").PasteSpecial(-4122)
$totals.Range("A7").Value = 5
$totals.Range("B7").Value = "2020-Q4"
$totals.Range("C7").Value = 1
$totals.Range("D7").Value = 0.44

# ------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet. Duplicate the existing
#    "2022-Q2" sheet (same column layout/styling) right before itself,
#    rename it, then overwrite the fund figures with the new quarter's
#    numbers.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("C2").Value = "华安国际龙头（DAX）ETF（QDII）"

$q3.Range("D2:G2").NumberFormat = "@"
$q3.Range("D2").Value = "5.54"
$q3.Range("E2").Value = "93.57"
$q3.Range("F2").Value = "4.31"
$q3.Range("G2").Value = "0.2388"

$q3.Range("H2").Value = 7
